$d = $word.ActiveDocument

# --- Change 1: "Steps to install Visual Studio Code." -> "Steps to install Spyder IDE."
# plus a _GoBack bookmark placed right after "Spyder IDE" (before the final ".")
# The first paragraph in the document holds this heading text.
$p1 = $d.Paragraphs(1)
$p1Start = $p1.Range.Start
$p1End = $p1.Range.End

$headingRange = $d.Range($p1Start, $p1End)
$found = $headingRange.Find.Execute("Steps to install Visual Studio Code.", $true, $false, $false, $false, $false, $true, 1, $false, "Steps to install Spyder IDE.", 2)

# Force a run split right before "Spyder IDE" using a temporary bookmark, then
# remove it once the real bookmark has been placed so only "_GoBack" remains.
$splitPos = $p1Start + 17
$goBackPos = $p1Start + 27

$tempRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("TempSplitMarker", $tempRange)

$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)

$d.Bookmarks("TempSplitMarker").Delete()

# --- Change 2: merge the split run-pair """ + "." into a single ""."  run
# in the first "Click on "Next"." paragraph (the one that still has them split).
$leftQuote = [char]0x201C
$rightQuote = [char]0x201D
$target = $rightQuote + "."

$wanted = "Click on ${leftQuote}Next${rightQuote}."
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd() -eq $wanted) {
        $pRange = $d.Range($p.Range.Start, $p.Range.End)
        $null = $pRange.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, $target, 2)
        break
    }
}
